$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$new = $wb.Worksheets.Add($null, $ws1)
$new.Name = "Sheet1"
$new.Range("Z1").HorizontalAlignment = -4152
$new.Range("A1:R8").HorizontalAlignment = -4108
